$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.009957909584045
$ws.Range("B1").Value = 2.130551099777222
$ws.Range("C1").Value = 7.224740028381348
$ws.Range("D1").Value = 2.381644248962402
$ws.Range("E1").Value = 1.339997291564941
